$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value  = "I'm Not The Only One"
$ws.Range("B8").Value  = "It Ain" + [char]0x2019 + "t Me (with Selena Gomez)"
$ws.Range("B12").Value = "XO Tour Llif3"
$ws.Range("B13").Value = "Ni**as In Paris"
$ws.Range("B34").Value = "Wild Thoughts (feat. Rihanna & Bryson Tiller)"
$ws.Range("B39").Value = "rockstar (feat. 21 Savage)"
$ws.Range("B41").Value = "2U (Originally Performed by David Guetta Feat. Justin Bieber) - Karaoke Version"
$ws.Range("B42").Value = "Ballin' (with Roddy Ricch)"
$ws.Range("B44").Value = "Feels Like Summer"
$ws.Range("B45").Value = "Mama Cry"
$ws.Range("B49").Value = "Chantaje (feat. Maluma)"
$ws.Range("B51").Value = "Havana (feat. Young Thug)"
$ws.Range("B64").Value = "You Don't Know Me (feat. Duane Harden) - Radio Edit"
$ws.Range("B68").Value = "How Far I'll Go (From " + [char]0x22 + "Moana" + [char]0x22 + ")"
$ws.Range("B85").Value = "side to side - live"
$ws.Range("B86").Value = "Otra vez (feat. J Balvin)"
$ws.Range("B90").Value = "Malibu Nights"
$ws.Range("B91").Value = "Up All Night"
$ws.Range("B92").Value = "Hear Me Now (feat. DIAMANTE)"
$ws.Range("B95").Value = "Friends (with BloodPop" + [char]0x00AE + ")"
